# Commit: "Add files via upload" — fills in the last 4 "1.9.x" task-card
# sheets (Entree(s)/Sortie(s)/Debut/Fin) that were left blank, and nudges
# the in-memory cursor/selection on a few sheets to match where the author
# last clicked while editing.

$wb = $excel.ActiveWorkbook

# --- 1.8.3: no data changed, user just left the selection on C7 ---
$ws183 = $wb.Worksheets.Item("1.8.3")
$ws183.Range("C7").Select()

# --- 1.9.1 ---
$ws191 = $wb.Worksheets.Item("1.9.1")
$ws191.Range("C7").Value = "Définition de l'innovation "
$ws191.Range("C6").Value = "Banc d'esssai et compétences de la session"
$ws191.Range("C8").Value = "06/08/2016"
$ws191.Range("C9").Value = "08/01/2016"
$ws191.Range("C7").Select()

# --- 1.9.2 ---
$ws192 = $wb.Worksheets.Item("1.9.2")
$ws192.Range("C6").Value = "Définition de l'innovation "
$ws192.Range("C7").Value = "Code/Simulation "
$ws192.Range("C8").Value = "06/08/2016"
$ws192.Range("C9").Value = "08/01/2016"
$ws192.Range("C6").Select()

# --- 1.9.3 ---
$ws193 = $wb.Worksheets.Item("1.9.3")
$ws193.Range("C6").Value = "Définition de l'innovation "
$ws193.Range("C7").Value = "Fonctionnement du banc d'essai"
$ws193.Range("C8").Value = "07/08/2016"
$ws193.Range("C9").Value = "08/01/2016"
$ws193.Range("C6").Select()

# --- 1.9.4 ---
$ws194 = $wb.Worksheets.Item("1.9.4")
$ws194.Range("C6").Value = "Définition de l'innovation "
$ws194.Range("C7").Value = "Visualisation sur Blender"
# This one picked up a heavier top border (it now visually closes off the
# thick-bottomed row above it), matching the author's original formatting.
$ws194.Range("C7").Borders.Item(8).Weight = -4138
$ws194.Range("C8").Value = "06/08/2016"
$ws194.Range("C9").Value = "08/01/2016"
$ws194.Range("H13").Select()

# NOTE: the diff also updates the x15ac:absPath breadcrumb (the folder the
# file lived in on the author's machine, "...\Projet\" -> "...\Projet\Gestion\").
# That value is stamped by the real Excel client from the OS save location and
# isn't exposed anywhere on the Workbook/Application object model (Path,
# FullName, SaveAs(...) are all inert no-ops against it here), so it can't be
# reproduced through COM automation.
